$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.828.21'
$ws.Cells.Item(2, 5).Value = '  -0.10%  '

$ws.Cells.Item(3, 4).Value = '2.248.16'
$ws.Cells.Item(3, 5).Value = '  +0.52%  '

$ws.Cells.Item(4, 5).Value = '  +0.17%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '112.57'
$ws.Cells.Item(5, 5).Value = '  -1.49%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '297.19'
$ws.Cells.Item(6, 5).Value = '  +7.32%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.629'
$ws.Cells.Item(7, 5).Value = '  +0.42%  '

$ws.Cells.Item(8, 5).Value = '  -0.29%  '

$ws.Cells.Item(9, 5).Value = '  -0.19%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '43.99'
$ws.Cells.Item(10, 5).Value = '  -5.37%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0926'
$ws.Cells.Item(11, 5).Value = '  +0.02%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '54.26'
$ws.Cells.Item(12, 5).Value = '  +0.66%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '8.98'
$ws.Cells.Item(13, 5).Value = '  -0.74%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '1.08'
$ws.Cells.Item(14, 5).Value = '  +23.25%  '

$ws.Cells.Item(15, 5).Value = '  -0.63%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '15.24'
$ws.Cells.Item(16, 5).Value = '  -0.13%  '

$ws.Cells.Item(17, 4).Value = '2.588.51'
$ws.Cells.Item(17, 5).Value = '  +0.52%  '

$ws.Cells.Item(18, 4).Value = '2.277.62'
$ws.Cells.Item(18, 5).Value = '  +1.49%  '

$ws.Cells.Item(19, 4).Value = '42.766.45'
$ws.Cells.Item(19, 5).Value = '  -0.13%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.0000107'
$ws.Cells.Item(20, 5).Value = '  -0.38%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '7.18'
$ws.Cells.Item(21, 5).Value = '  +6.34%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '74.81'
$ws.Cells.Item(22, 5).Value = '  +3.74%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.52'
$ws.Cells.Item(23, 5).Value = '  +17.95%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.46'
$ws.Cells.Item(24, 5).Value = '  +5.06%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '252.02'
$ws.Cells.Item(25, 5).Value = '  +9.01%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '8.95'
$ws.Cells.Item(26, 5).Value = '  -3.09%  '

$ws.Cells.Item(27, 5).Value = '  -0.47%  '

$ws.Cells.Item(28, 5).Value = '  -3.17%  '

$ws.Cells.Item(29, 5).Value = '  -0.59%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '37.97'
$ws.Cells.Item(30, 5).Value = '  -5.45%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '175.37'
$ws.Cells.Item(31, 5).Value = '  +1.24%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '22.12'
$ws.Cells.Item(32, 5).Value = '  +5.05%  '

$ws.Cells.Item(33, 5).Value = '  -3.36%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.0890'
$ws.Cells.Item(34, 5).Value = '  -0.10%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '5.67'
$ws.Cells.Item(35, 5).Value = '  +1.97%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '5.08'
$ws.Cells.Item(36, 5).Value = '  +9.29%  '

$ws.Cells.Item(37, 5).Value = '  +0.07%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '4.24'
$ws.Cells.Item(38, 5).Value = '  -4.05%  '

$ws.Cells.Item(39, 5).Value = '  +1.44%  '

$ws.Cells.Item(40, 5).Value = '  -1.48%  '

$ws.Cells.Item(41, 5).Value = '  -5.34%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '72.10'
$ws.Cells.Item(42, 5).Value = '  +1.64%  '

$ws.Cells.Item(43, 5).Value = '  -0.53%  '

$ws.Cells.Item(44, 5).Value = '  +0.05%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '12.49'
$ws.Cells.Item(45, 5).Value = '  -5.17%  '

$ws.Cells.Item(46, 5).Value = '  -0.64%  '

$ws.Cells.Item(47, 5).Value = '  -2.25%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '106.46'
$ws.Cells.Item(48, 5).Value = '  +6.17%  '

$ws.Cells.Item(49, 5).Value = '  +2.55%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '8.66'
$ws.Cells.Item(50, 5).Value = '  +2.64%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '72.04'
$ws.Cells.Item(51, 5).Value = '  +2.87%  '
